$wb = $excel.ActiveWorkbook

# --- Sheet "Euramet": update raw acquisition data for rows 7 and 8 ---
$ws1 = $wb.Worksheets.Item("Euramet")

$ws1.Range("E7").Value = 989
$ws1.Range("F7").Value = -0.3
$ws1.Range("G7").Value = -2.4
$ws1.Range("H7").Value = 2.42215

$ws1.Range("E8").Value = 883
$ws1.Range("F8").Value = 257
$ws1.Range("G8").Value = -253.9
$ws1.Range("H8").Value = 1.90669

# Work on the executable/acquisition paused: the remaining measurement
# rows (9-19 and 29-41) had not yet been re-acquired, so their raw
# D:H values are cleared back out (formatting/style stays untouched).
$ws1.Range("D9:H19").ClearContents()
$ws1.Range("D29:H41").ClearContents()

# --- Sheet "Istruzioni Uso": placeholder values reset to "-" ---
$ws2 = $wb.Worksheets.Item("Istruzioni Uso")
$ws2.Range("B65:B69").Value = "-"
